# edit.ps1
# Applies the DTR date-range shift (03-21..04-03-2015 -> 04-25..05-08-2015)
# plus related formula/style fixes, as described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Re-style rows 7,8,9,10,11,15,16,17,18 BEFORE we touch their values,
#    using Copy + PasteSpecial(xlPasteFormats) from donor rows that
#    already carry the exact target style index, so Excel reuses the
#    existing cellXfs entry instead of allocating a new one.
#    xlPasteFormats = -4122
# ---------------------------------------------------------------------
$xlPasteFormats = -4122

# Donor for style "9" (red/holiday look): row 17 (A17:P17) -- grab this
# BEFORE row 17 itself is restyled below.
$ws.Range("A17:P17").Copy()
$ws.Range("A7:P7").PasteSpecial($xlPasteFormats)
$ws.Range("A8:P8").PasteSpecial($xlPasteFormats)
$ws.Range("A11:P11").PasteSpecial($xlPasteFormats)

# Donor for plain style "6": row 5 (A5:P5) stays style 6 throughout.
$ws.Range("A5:P5").Copy()
$ws.Range("A9:P9").PasteSpecial($xlPasteFormats)
$ws.Range("A15:P15").PasteSpecial($xlPasteFormats)
$ws.Range("A16:P16").PasteSpecial($xlPasteFormats)
$ws.Range("A17:P17").PasteSpecial($xlPasteFormats)
$ws.Range("A18:P18").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Helper: write a literal text value into a cell without Excel's
#    automatic date-pattern recognition turning "MM-DD-YYYY" strings
#    into date serials. We briefly force Text format, assign, then
#    restore the exact original style via a format-only paste from a
#    same-styled neighbour cell.
# ---------------------------------------------------------------------
function Set-TextValue($rangeAddr, $text, $styleDonorAddr) {
    $cell = $ws.Range($rangeAddr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $ws.Range($styleDonorAddr).Copy()
    $cell.PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------
# 3) Row 7 (04-27-2015 / Monday / VL trip, style 9)
# ---------------------------------------------------------------------
Set-TextValue "A7" "04-27-2015" "B7"
$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("H7").Value = 1
Set-TextValue "P7" "~VL Trip to Baguio. :)" "B7"

# ---------------------------------------------------------------------
# 4) Row 8 (04-28-2015 / Tuesday / VL trip, style 9)
# ---------------------------------------------------------------------
Set-TextValue "A8" "04-28-2015" "B8"
$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("H8").Value = 1
Set-TextValue "P8" "~VL Trip to Baguio. :)" "B8"

# ---------------------------------------------------------------------
# 5) Row 9 (04-29-2015 / Wednesday, style 6)
# ---------------------------------------------------------------------
Set-TextValue "A9" "04-29-2015" "B9"
Set-TextValue "C9" "08:02:45" "B9"
Set-TextValue "D9" "18:47:37" "B9"

# ---------------------------------------------------------------------
# 6) Row 10 (04-30-2015 / Thursday, style 6 - unchanged)
# ---------------------------------------------------------------------
Set-TextValue "A10" "04-30-2015" "B10"
Set-TextValue "C10" "08:57:09" "B10"
Set-TextValue "D10" "18:43:31" "B10"
$ws.Range("E10").Value = 0.5

# ---------------------------------------------------------------------
# 7) Row 11 (05-01-2015 / Friday / Labor Day holiday, style 9)
# ---------------------------------------------------------------------
Set-TextValue "A11" "05-01-2015" "B11"
$ws.Range("C11").ClearContents()
$ws.Range("D11").ClearContents()
Set-TextValue "P11" "~Special Day - Labor ( Regular Holiday )" "B11"

# ---------------------------------------------------------------------
# 8) Row 12 (05-02-2015 / Saturday, style 6, only date changes)
# ---------------------------------------------------------------------
Set-TextValue "A12" "05-02-2015" "B12"

# ---------------------------------------------------------------------
# 9) Row 13 (05-03-2015 / Sunday, style 6, only date changes)
# ---------------------------------------------------------------------
Set-TextValue "A13" "05-03-2015" "B13"

# ---------------------------------------------------------------------
# 10) Row 14 (05-04-2015 / Monday, style 6)
# ---------------------------------------------------------------------
Set-TextValue "A14" "05-04-2015" "B14"
Set-TextValue "C14" "08:25:44" "B14"
Set-TextValue "D14" "19:01:02" "B14"

# ---------------------------------------------------------------------
# 11) Row 15 (05-05-2015 / Tuesday, style 6 - newly converted from holiday row)
# ---------------------------------------------------------------------
Set-TextValue "A15" "05-05-2015" "B15"
Set-TextValue "C15" "08:41:28" "B15"
Set-TextValue "D15" "18:45:38" "B15"
$ws.Range("E15").Value = 0.25

# ---------------------------------------------------------------------
# 12) Row 16 (05-06-2015 / Wednesday, style 6)
# ---------------------------------------------------------------------
Set-TextValue "A16" "05-06-2015" "B16"
Set-TextValue "C16" "08:26:35" "B16"
Set-TextValue "D16" "18:36:11" "B16"

# ---------------------------------------------------------------------
# 13) Row 17 (05-07-2015 / Thursday, style 6 - no more holiday remark)
# ---------------------------------------------------------------------
Set-TextValue "A17" "05-07-2015" "B17"
Set-TextValue "C17" "08:25:05" "B17"
Set-TextValue "D17" "18:40:35" "B17"
$ws.Range("P17").ClearContents()

# ---------------------------------------------------------------------
# 14) Row 18 (05-08-2015 / Friday, style 6 - no more holiday remark)
# ---------------------------------------------------------------------
Set-TextValue "A18" "05-08-2015" "B18"
Set-TextValue "C18" "08:15:12" "B18"
Set-TextValue "D18" "17:34:26" "B18"
$ws.Range("P18").ClearContents()

# ---------------------------------------------------------------------
# 15) Formula range updates (E5:E15 -> E5:E18, etc.)
# ---------------------------------------------------------------------
$ws.Range("E19").Formula = "=COUNT(E5:E18)"
$ws.Range("E20").Formula = "=SUM(E5:E18)"
$ws.Range("H22").Formula = "=SUM(H5:H18)"
$ws.Range("I22").Formula = "=SUM(I5:I18)"

# I22 also regains the formula-result style (s=17) matching H22/G21/etc.
$ws.Range("H22").Copy()
$ws.Range("I22").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 16) I23: accumulated-OT manual value 2.5 -> 0
# ---------------------------------------------------------------------
$ws.Range("I23").Value = 0

# ---------------------------------------------------------------------
# 17) C26 formula: FLOOR's 3rd "significance" arg removed (2-arg FLOOR)
# ---------------------------------------------------------------------
$ws.Range("C26").Formula = '=FLOOR(I23,1)&"."&FLOOR(MOD(I23*8,8),1)&"."&(MOD(I23*8,8)-FLOOR(MOD(I23*8,8),1))*60'
